$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (same shared-string slots, new Arabic labels):
#   A1: "id"      -> "الاسم"   (Name)
#   C1: "address" -> "العنوان" (Address)
#   F1: "record"  -> "السجل"   (Record)
$ws.Range("A1").Value = "الاسم"
$ws.Range("C1").Value = "العنوان"
$ws.Range("F1").Value = "السجل"

# Scroll the view so column C is the left-most visible column
# (sheetView topLeftCell goes from D1 to C1).
$excel.ActiveWindow.ScrollColumn = 3
